# Simulator full-month coverage, persist logs, fix employees
#
# Re-generates the 2026-01-19 week for Jason Green: the former PTO day
# becomes a full Regular day, client assignments are refreshed, the
# weekly hours roll up to a full 40 (was 32), and the employee id is
# reassigned.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")
$ws2 = $wb.Worksheets.Item("Jason Schema")

# ---------------------------------------------------------------------
# "Weekly Timesheet" sheet: Date | Client | Hours | Type | Rate | Total
# ---------------------------------------------------------------------

# Row 2 - 2026-01-19 (was PTO, now a full Regular day)
$ws1.Range("B2").Value = "Regular"
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = "Regular"
$ws1.Range("E2").Value = 100
$ws1.Range("F2").Value = 800

# Row 3 - 2026-01-20
$ws1.Range("B3").Value = "Leixner/Smith"
$ws1.Range("C3").Value = 8
$ws1.Range("D3").Value = "Regular"
$ws1.Range("E3").Value = 100
$ws1.Range("F3").Value = 800

# Row 4 - 2026-01-21
$ws1.Range("B4").Value = "Hunter"
$ws1.Range("C4").Value = 8
$ws1.Range("D4").Value = "Regular"
$ws1.Range("E4").Value = 100
$ws1.Range("F4").Value = 800

# Row 5 - 2026-01-22
$ws1.Range("B5").Value = "O'Connor"
$ws1.Range("C5").Value = 8
$ws1.Range("D5").Value = "Regular"
$ws1.Range("E5").Value = 100
$ws1.Range("F5").Value = 800

# Row 6 - 2026-01-23
$ws1.Range("B6").Value = "Varricchio"
$ws1.Range("C6").Value = 8
$ws1.Range("D6").Value = "Regular"
$ws1.Range("E6").Value = 100
$ws1.Range("F6").Value = 800

# Row 8 - SUBTOTAL: hours 32 -> 40, amount 0 -> 4000 (note text follows)
$ws1.Range("C8").Value = 40
$ws1.Range("D8").Value = "Reg: 40 / OT: 0"
$ws1.Range("F8").Value = 4000

# Row 11 - HOURLY SUBTOTAL amount
$ws1.Range("F11").Value = 4000

# Row 13 - GRAND TOTAL amount
$ws1.Range("F13").Value = 4000

# ---------------------------------------------------------------------
# "Jason Schema" sheet:
# Employee | Employee ID | Date | Client | Hours | Rate | Total | Type | Notes
# ---------------------------------------------------------------------

# Employee ID reassigned for every data row
$ws2.Range("B2:B6").Value = "emp_qhpjptqm"

# Row 2 - 2026-01-19 (was PTO, now a full Regular day)
$ws2.Range("D2").Value = "Markfield"
$ws2.Range("E2").Value = 8
$ws2.Range("F2").Value = 100
$ws2.Range("G2").Value = 800
$ws2.Range("H2").Value = "Regular"
$ws2.Range("I2").Value = ""

# Row 3 - 2026-01-20
$ws2.Range("D3").Value = "Leixner/Smith"
$ws2.Range("E3").Value = 8
$ws2.Range("F3").Value = 100
$ws2.Range("G3").Value = 800
$ws2.Range("H3").Value = "Regular"

# Row 4 - 2026-01-21
$ws2.Range("D4").Value = "Hunter"
$ws2.Range("E4").Value = 8
$ws2.Range("F4").Value = 100
$ws2.Range("G4").Value = 800
$ws2.Range("H4").Value = "Regular"

# Row 5 - 2026-01-22
$ws2.Range("D5").Value = "O'Connor"
$ws2.Range("E5").Value = 8
$ws2.Range("F5").Value = 100
$ws2.Range("G5").Value = 800
$ws2.Range("H5").Value = "Regular"

# Row 6 - 2026-01-23
$ws2.Range("D6").Value = "Varricchio"
$ws2.Range("E6").Value = 8
$ws2.Range("F6").Value = 100
$ws2.Range("G6").Value = 800
$ws2.Range("H6").Value = "Regular"
